# Generate Report for Handback
# - Flip the Overview "Handed back" status from "in sync" to "not in sync"
#   with en-US (this text is duplicated across the Overview sheet and the
#   per-locale "Status" column on the zh-cn / de-de sheets).
# - Refresh the "Correspond Handback DateTime" for the newest handback on
#   the zh-cn and de-de sheets (row 3 = the 66b333c0... file).
# - Widen the Status columns to fit the longer status text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("E3").Value = $newStatus
$ws.Range("F3").Value = $newStatus

# Column widths grew to fit the longer status text (~33.46 chars); the
# closest width this engine's column-width quantization can reach.
$ws.Columns.Item(5).ColumnWidth = 32.666666666666664
$ws.Columns.Item(6).ColumnWidth = 32.666666666666664

# ---- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus
$ws.Range("K3").Value = "2016-09-06 11:27:09"
$ws.Columns.Item(3).ColumnWidth = 32.666666666666664

# ---- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus
$ws.Range("K3").Value = "2016-09-06 11:27:32"
$ws.Columns.Item(3).ColumnWidth = 32.666666666666664
